$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# row 4
$ws.Range("H4").Value = 124
$ws.Range("I4").Value = 131.85715
$ws.Range("K4").Value = 131.85715
$ws.Range("M4").Value = -17.85714999999999

# row 43
$ws.Range("H43").Value = 1199.3334

# row 52
$ws.Range("H52").Value = 2199.4
$ws.Range("I52").Value = 998
$ws.Range("K52").Value = 2994
$ws.Range("M52").Value = -2834

# row 103
$ws.Range("H103").Value = 1880
$ws.Range("I103").Value = 1900
$ws.Range("K103").Value = 5700
$ws.Range("M103").Value = -5114

# row 137
$ws.Range("H137").Value = 2612.9333
$ws.Range("I137").Value = 2513.3572
$ws.Range("K137").Value = 7540.071599999999
$ws.Range("M137").Value = -4990.071599999999

# row 138
$ws.Range("H138").Value = 2775.842
$ws.Range("I138").Value = 1581
$ws.Range("K138").Value = 4743
$ws.Range("M138").Value = 397

$ws = $wb.Worksheets.Item("ARM")
# row 32
$ws.Range("H32").Value = 2457.1404
$ws.Range("I32").Value = 2273.7637
$ws.Range("K32").Value = 2273.7637
$ws.Range("M32").Value = -1986.7637

# row 60
$ws.Range("H60").Value = 0
$ws.Range("I60").Value = 0
$ws.Range("K60").Value = 0
$ws.Range("M60").Value = ""

# row 63
$ws.Range("H63").Value = 2860
$ws.Range("I63").Value = 2600
$ws.Range("J63").Value = 3250
$ws.Range("K63").Value = 2600
$ws.Range("L63").Value = 3250
$ws.Range("M63").Value = -1914
$ws.Range("N63").Value = -4622

# row 66
$ws.Range("H66").Value = 2860
$ws.Range("I66").Value = 2600
$ws.Range("J66").Value = 3250
$ws.Range("K66").Value = 13000
$ws.Range("L66").Value = 16250
$ws.Range("M66").Value = -9568
$ws.Range("N66").Value = -23114

# row 70
$ws.Range("H70").Value = 0
$ws.Range("J70").Value = 0
$ws.Range("L70").Value = ""
$ws.Range("N70").Value = ""

# row 73
$ws.Range("H73").Value = 0
$ws.Range("J73").Value = 0
$ws.Range("L73").Value = ""
$ws.Range("N73").Value = ""

$ws = $wb.Worksheets.Item("BSM")
# row 20
$ws.Range("H20").Value = 3942.1428
$ws.Range("I20").Value = 3919.4
$ws.Range("J20").Value = 3999
$ws.Range("K20").Value = 3919.4
$ws.Range("L20").Value = 3999
$ws.Range("M20").Value = -3672.4
$ws.Range("N20").Value = -4493

# row 82
$ws.Range("H82").Value = 31279.666
$ws.Range("J82").Value = 39661.11
$ws.Range("L82").Value = 39661.11
$ws.Range("N82").Value = -40427.11

# row 85
$ws.Range("H85").Value = 31279.666
$ws.Range("J85").Value = 39661.11
$ws.Range("L85").Value = 39661.11
$ws.Range("N85").Value = -42313.11

# row 92
$ws.Range("H92").Value = 30000
$ws.Range("J92").Value = 30000
$ws.Range("L92").Value = 30000
$ws.Range("N92").Value = -34992

# row 134
$ws.Range("H134").Value = 13111.143
$ws.Range("I134").Value = 13111.143
$ws.Range("K134").Value = 39333.429
$ws.Range("M134").Value = -36798.429

$ws = $wb.Worksheets.Item("CRP")
# row 31
$ws.Range("H31").Value = 1803
$ws.Range("I31").Value = 1344.5
$ws.Range("J31").Value = 6388
$ws.Range("K31").Value = 1344.5
$ws.Range("L31").Value = 6388
$ws.Range("M31").Value = -1049.5
$ws.Range("N31").Value = -6978

# row 34
$ws.Range("H34").Value = 1803
$ws.Range("I34").Value = 1344.5
$ws.Range("J34").Value = 6388
$ws.Range("K34").Value = 1344.5
$ws.Range("L34").Value = 6388
$ws.Range("M34").Value = -1142.5
$ws.Range("N34").Value = -6792

# row 58
$ws.Range("H58").Value = 2843.4707
$ws.Range("I58").Value = 2042.9231
$ws.Range("J58").Value = 5445.25
$ws.Range("K58").Value = 2042.9231
$ws.Range("L58").Value = 5445.25
$ws.Range("M58").Value = -1839.9231
$ws.Range("N58").Value = -5851.25

# row 105
$ws.Range("H105").Value = 366
$ws.Range("I105").Value = 366
$ws.Range("K105").Value = 366
$ws.Range("M105").Value = 1381

# row 136
$ws.Range("H136").Value = 2843.4707
$ws.Range("I136").Value = 2042.9231
$ws.Range("J136").Value = 5445.25
$ws.Range("K136").Value = 6128.7693
$ws.Range("L136").Value = 16335.75
$ws.Range("M136").Value = -3578.7693
$ws.Range("N136").Value = -21435.75

$ws = $wb.Worksheets.Item("GSM")
# row 58
$ws.Range("H58").Value = 20000
$ws.Range("J58").Value = 20000
$ws.Range("L58").Value = 20000
$ws.Range("N58").Value = -20554

# row 70
$ws.Range("H70").Value = 41673664
$ws.Range("I70").Value = 55561972
$ws.Range("K70").Value = 55561972
$ws.Range("M70").Value = -55561702

# row 73
$ws.Range("H73").Value = 41673664
$ws.Range("I73").Value = 55561972
$ws.Range("K73").Value = 55561972
$ws.Range("M73").Value = -55561036

# row 80
$ws.Range("H80").Value = 2397.8
$ws.Range("I80").Value = 1700
$ws.Range("J80").Value = 2572.25
$ws.Range("K80").Value = 1700
$ws.Range("L80").Value = 2572.25
$ws.Range("M80").Value = -702
$ws.Range("N80").Value = -4568.25

# row 83
$ws.Range("H83").Value = 2397.8
$ws.Range("I83").Value = 1700
$ws.Range("J83").Value = 2572.25
$ws.Range("K83").Value = 8500
$ws.Range("L83").Value = 12861.25
$ws.Range("M83").Value = -3508
$ws.Range("N83").Value = -22845.25

# row 97
$ws.Range("H97").Value = 820.1667
$ws.Range("I97").Value = 806.125
$ws.Range("J97").Value = 848.25
$ws.Range("K97").Value = 806.125
$ws.Range("L97").Value = 848.25
$ws.Range("M97").Value = -310.125
$ws.Range("N97").Value = -1840.25

# row 132
$ws.Range("H132").Value = 3587.2104
$ws.Range("I132").Value = 3587.2104
$ws.Range("K132").Value = 10761.6312
$ws.Range("M132").Value = -8231.6312

# row 134
$ws.Range("H134").Value = 97499.5
$ws.Range("J134").Value = 97499.5
$ws.Range("L134").Value = 292498.5
$ws.Range("N134").Value = -297568.5

$ws = $wb.Worksheets.Item("LTW")
# row 132
$ws.Range("H132").Value = 3095.7222
$ws.Range("I132").Value = 2517.1538
$ws.Range("J132").Value = 4600
$ws.Range("K132").Value = 7551.4614
$ws.Range("L132").Value = 13800
$ws.Range("M132").Value = -5021.4614
$ws.Range("N132").Value = -18860

$ws = $wb.Worksheets.Item("WVR")
# row 93
$ws.Range("H93").Value = 39000
$ws.Range("J93").Value = 39000
$ws.Range("L93").Value = 39000
$ws.Range("N93").Value = -43992

# row 116
$ws.Range("H116").Value = 0
$ws.Range("J116").Value = 0
$ws.Range("L116").Value = ""
$ws.Range("N116").Value = ""

# row 117
$ws.Range("H117").Value = 21500
$ws.Range("J117").Value = 21500
$ws.Range("L117").Value = 21500
$ws.Range("N117").Value = -30678

# row 136
$ws.Range("H136").Value = 3788.15
$ws.Range("I136").Value = 3711.2727
$ws.Range("K136").Value = 11133.8181
$ws.Range("M136").Value = -8583.8181
